# Indicators_Targets_ext.xlsx - header renamed to lowercase/underscore naming
# convention for consistency, and selection updated to the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row (row 1) from CamelCase to lower_snake_case names.
$ws.Range("A1").Value = "indicator_name"
$ws.Range("B1").Value = "actual"
$ws.Range("C1").Value = "actual_lastweek"
$ws.Range("D1").Value = "actual_lastyear"
$ws.Range("E1").Value = "target"

# Update the selected range to the header row.
$ws.Range("A1:E1").Select()
